$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.490.31"
$ws.Range('E2').Value = "'  +1.34%  "
$ws.Range('D3').Value = "'3.398.10"
$ws.Range('E3').Value = "'  +3.89%  "
$ws.Range('E4').Value = "'  -0.16%  "
$ws.Range('D5').Value = "'575.43"
$ws.Range('E5').Value = "'  +2.38%  "
$ws.Range('D6').Value = "'138.37"
$ws.Range('E6').Value = "'  +9.41%  "
$ws.Range('E7').Value = "'  -0.07%  "
$ws.Range('D8').Value = "'3.396.81"
$ws.Range('E8').Value = "'  +3.80%  "
$ws.Range('D9').Value = "'0.477"
$ws.Range('E9').Value = "'  +0.72%  "
$ws.Range('D10').Value = "'7.60"
$ws.Range('E10').Value = "'  +4.40%  "
$ws.Range('E11').Value = "'  +8.27%  "
$ws.Range('E12').Value = "'  +6.44%  "
$ws.Range('D13').Value = "'3.983.78"
$ws.Range('E13').Value = "'  +3.67%  "
$ws.Range('E14').Value = "'  +1.72%  "
$ws.Range('D15').Value = "'0.0000180"
$ws.Range('E15').Value = "'  +8.12%  "
$ws.Range('D16').Value = "'3.404.43"
$ws.Range('E16').Value = "'  +3.74%  "
$ws.Range('D17').Value = "'25.46"
$ws.Range('E17').Value = "'  +5.36%  "
$ws.Range('D18').Value = "'61.586.00"
$ws.Range('E18').Value = "'  +1.15%  "
$ws.Range('D19').Value = "'14.11"
$ws.Range('E19').Value = "'  +7.04%  "
$ws.Range('D20').Value = "'5.90"
$ws.Range('E20').Value = "'  +5.34%  "
$ws.Range('E21').Value = "'  +4.45%  "
$ws.Range('D22').Value = "'390.13"
$ws.Range('E22').Value = "'  +11.94%  "
$ws.Range('E23').Value = "'  +4.13%  "
$ws.Range('D24').Value = "'3.538.94"
$ws.Range('E24').Value = "'  +3.72%  "
$ws.Range('E25').Value = "'  +0.21%  "
$ws.Range('D26').Value = "'0.0000127"
$ws.Range('E26').Value = "'  +18.62%  "
$ws.Range('D27').Value = "'71.04"
$ws.Range('E27').Value = "'  +2.68%  "
$ws.Range('E28').Value = "'  +13.68%  "
$ws.Range('D29').Value = "'7.74"
$ws.Range('E29').Value = "'  +8.29%  "
$ws.Range('E30').Value = "'  +0.36%  "
$ws.Range('D31').Value = "'8.34"
$ws.Range('E31').Value = "'  +7.36%  "
$ws.Range('D32').Value = "'0.162"
$ws.Range('E32').Value = "'  +8.95%  "
$ws.Range('E33').Value = "'  +3.19%  "
$ws.Range('E34').Value = "'  -0.09%  "
$ws.Range('D35').Value = "'3.430.45"
$ws.Range('E35').Value = "'  +3.67%  "
$ws.Range('D36').Value = "'23.57"
$ws.Range('E36').Value = "'  +4.68%  "
$ws.Range('D37').Value = "'5.52"
$ws.Range('E37').Value = "'  +5.98%  "
$ws.Range('D38').Value = "'7.00"
$ws.Range('E38').Value = "'  +3.51%  "
$ws.Range('E39').Value = "'  +5.53%  "
$ws.Range('D40').Value = "'162.01"
$ws.Range('E40').Value = "'  +1.52%  "
$ws.Range('D41').Value = "'0.0793"
$ws.Range('E41').Value = "'  +6.16%  "
$ws.Range('B42').Value = "'FirstDigitalUSD"
$ws.Range('C42').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = "'  -0.19%  "
$ws.Range('B43').Value = "'Stacks"
$ws.Range('C43').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('D43').Value = "'1.73"
$ws.Range('E43').Value = "'  +12.10%  "
$ws.Range('B44').Value = "'ONDO"
$ws.Range('C44').Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range('D44').Value = "'1.23"
$ws.Range('E44').Value = "'  +9.29%  "
$ws.Range('D45').Value = "'4.48"
$ws.Range('E45').Value = "'  +4.46%  "
$ws.Range('D46').Value = "'0.774"
$ws.Range('E46').Value = "'  +5.16%  "
$ws.Range('D47').Value = "'41.14"
$ws.Range('E47').Value = "'  +0.43%  "
$ws.Range('D48').Value = "'24.21"
$ws.Range('E48').Value = "'  +8.09%  "
$ws.Range('E49').Value = "'  +5.03%  "
$ws.Range('D50').Value = "'23.00"
$ws.Range('E50').Value = "'  +8.48%  "
$ws.Range('D51').Value = "'2.379.71"
$ws.Range('E51').Value = "'  +10.75%  "
